$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook had three separate rows for the "Bonaire, Sint Eustatius and
# Saba" special municipality (Bonaire / Sint Eustatius / Saba), each using
# placeholder alpha_3 codes BESB/BESE/BESS. The Sint Eustatius row (row 8)
# is being dropped while work on the Saba/Eustatia split continues, and the
# remaining two rows get corrected alpha_3 codes (BON / ESS).

# Delete the "Sint Eustatius" row entirely (row 8) - this shifts every row
# below it up by one.
$ws.Rows.Item(8).Delete()

# Fix the alpha_3 codes for the two rows that remain (Bonaire on row 7, and
# Saba which has now shifted up into row 8).
$ws.Range("C7").Value = "BON"
$ws.Range("C8").Value = "ESS"

# Update the active selection left behind by the edit.
[void]$ws.Range("C8").Select()

# Turn on iterative calculation (used while resolving the Saba/Eustatia
# circular-reference issue mentioned in the commit message).
$wb.IterativeCalculation = $true
$wb.MaxIterations = 10
